$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Column1"
$ws.Range("E2").Value = "This should change to number/factors"

$ws.Range("E2").Select()
$ws.Application.ActiveWindow.ScrollColumn = 2
